# "Add milestone3 and update gantt"
#
# The author reworked the task schedule on the "Data" sheet:
#   - Front-End Development (row 7) now needs 20 days instead of 10
#   - Debugging Process (row 8) now needs 10 days instead of 6
#   - Creating Relational User Tables (row 9) now starts 45209 (was 45201)
#     and needs 8 days instead of 14
# The WORKDAY()/subtraction formulas in columns D and E, and the linked
# Gantt chart on the "Gantt Chart" sheet, pick up the new figures
# automatically on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

$ws.Range("C7").Value = 20
$ws.Range("C8").Value = 10
$ws.Range("B9").Value = 45209
$ws.Range("C9").Value = 8

$ws.Range("C9").Select() | Out-Null
